$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-0.2062081610468891, 0.4981741733236962, 0.4649261630253289, 0.6818549428033274, 0.657108150032778, 46)
    3  = @(-0.09587357185906692, 0.5000820061429816, 0.4618601616711875, 0.6796029441307531, 0.6804089494737995, 45)
    4  = @(-0.1372403718692351, 0.4340014692588762, 0.3663092181924378, 0.6052348454876321, 0.5962844087863356, 44)
    5  = @(-0.08352607004722953, 0.471891385230405, 0.4001742410426709, 0.6325932666750974, 0.6344757561618594, 43)
    6  = @(-0.1513393225090566, 0.4428429211668459, 0.3795138264528281, 0.616046935267783, 0.6044071860487836, 42)
    7  = @(-0.1148246166180287, 0.4704721399279175, 0.4088047900674556, 0.6393784404149514, 0.6367971418930269, 41)
    8  = @(-0.1576462645105448, 0.4410790206399262, 0.3843178504030022, 0.6199337467850917, 0.6071923665011824, 40)
    9  = @(-0.1265473388711615, 0.4734734975690142, 0.415038962515679, 0.6442351764035235, 0.6399417039254013, 39)
    10 = @(-0.181494460453781, 0.4408424176267598, 0.3840813816712615, 0.6197429964681017, 0.6005259724977672, 38)
    11 = @(-0.1278002099571802, 0.4693678920889392, 0.4189108628417125, 0.6472332368178511, 0.6432423264545254, 37)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
